# Update odds data in the active worksheet (Jogos_da_Semana_FlashScore_2025-03-20)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Union de Santa Fe vs Racing Club)
$ws.Range("G2").Value = 3.25
$ws.Range("I2").Value = 2.45
$ws.Range("J2").Value = 4
$ws.Range("L2").Value = 3.2
$ws.Range("M2").Value = 1.07
$ws.Range("O2").Value = 1.41
$ws.Range("P2").Value = 2.62
$ws.Range("Q2").Value = 1.74
$ws.Range("R2").Value = 1.99
$ws.Range("X2").Value = 1.15
$ws.Range("AD2").Value = 15
$ws.Range("AF2").Value = 34
$ws.Range("AJ2").Value = 6
$ws.Range("AN2").Value = 6.5
$ws.Range("AQ2").Value = 23
$ws.Range("AR2").Value = 21
$ws.Range("AS2").Value = 34

# Row 4 (Montevideo City vs Juventud)
$ws.Range("G4").Value = 1.91
$ws.Range("M4").Value = 1.07
$ws.Range("O4").Value = 1.4
$ws.Range("X4").Value = 1.22

# Row 5 (Penarol vs Liverpool M.)
$ws.Range("G5").Value = 1.5
$ws.Range("M5").Value = 1.07
$ws.Range("O5").Value = 1.33
$ws.Range("X5").Value = 1.25
